$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: "撰写第二轮迭代计划和第二轮评估报告。" ->
#         "撰写第二轮迭代计划和第二轮评估报告，通过第二轮迭代评审。"
# ------------------------------------------------------------------
$ok1 = $d.Content.Find.Execute(
    "撰写第二轮迭代计划和第二轮评估报告。", $false, $false, $false, $false,
    $false, $true, 1, $false,
    "撰写第二轮迭代计划和第二轮评估报告，通过第二轮迭代评审。", 2)
Write-Output "edit1: $ok1"

# ------------------------------------------------------------------
# Edit 2: append a sentence about the temperature sensor generator to
# the paragraph that ends with "...中实时折线图和热力图的组件。"
# ------------------------------------------------------------------
$ok2 = $d.Content.Find.Execute(
    "中实时折线图和热力图的组件。", $false, $false, $false, $false,
    $false, $true, 1, $false,
    "中实时折线图和热力图的组件。完成温度传感器的数据生成器和后端接口。", 2)
Write-Output "edit2: $ok2"

# ------------------------------------------------------------------
# Edit 3: insert a new paragraph describing the test-case / test-report
# work, sandwiched between two blank paragraphs, right after the
# (already blank) paragraph that follows
# "建立自动发送邮箱验证邮件的机制。".
# ------------------------------------------------------------------
$anchorIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "建立自动发送邮箱验证邮件的机制。`r") {
        $anchorIdx = $i
        break
    }
}
Write-Output "edit3 anchorIdx: $anchorIdx"

$emptyAfterAnchor = $d.Paragraphs.Item($anchorIdx + 1)
$emptyAfterAnchor.Range.InsertParagraphAfter()
$d.Paragraphs.Item($anchorIdx + 2).Range.Text = "记录登陆、注册、用户信息查看及修改业务的测试用例并撰写测试报告，含：单元测试、集成测试、安全性测试。"
$d.Paragraphs.Item($anchorIdx + 2).Range.InsertParagraphAfter()
Write-Output "edit3 paragraphs: $($d.Paragraphs.Count)"

# ------------------------------------------------------------------
# Edit 4: append a new trailing paragraph about the HTTP interceptor
# and logging, after the final content paragraph (the previously
# empty last paragraph becomes this new paragraph).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$last.Range.Text = "添加了HTTP请求的拦截器和日志，记录访问情况。"
Write-Output "edit4 paragraphs: $($d.Paragraphs.Count)"
